# "improve handle of empty multiform file"
# Add a new, empty "2013" form/profile sheet after the existing "2012" one,
# mirroring the same minimal single-cell "no instructions" placeholder the
# other empty profile sheets ("2012") use, and make it the active sheet.

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("2012")

# New sheet goes right after "2012" and becomes the active tab, exactly like
# Excel's normal "insert sheet after" behaviour.
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $source)
$newSheet.Name = "2013"

# Same placeholder text as the "2012" sheet's A1 (reuses the existing shared
# string "aucune instruction" rather than creating a new unique one).
$newSheet.Range("A1").Value = $source.Range("A1").Text

# Carry over the page setup from the sibling sheet so the new, empty sheet
# prints the same way the others do.
$newSheet.PageSetup.LeftMargin = $source.PageSetup.LeftMargin
$newSheet.PageSetup.RightMargin = $source.PageSetup.RightMargin
$newSheet.PageSetup.TopMargin = $source.PageSetup.TopMargin
$newSheet.PageSetup.BottomMargin = $source.PageSetup.BottomMargin
$newSheet.PageSetup.HeaderMargin = $source.PageSetup.HeaderMargin
$newSheet.PageSetup.FooterMargin = $source.PageSetup.FooterMargin
$newSheet.PageSetup.CenterHeader = $source.PageSetup.CenterHeader
$newSheet.PageSetup.CenterFooter = $source.PageSetup.CenterFooter
$newSheet.PageSetup.PaperSize = $source.PageSetup.PaperSize
$newSheet.PageSetup.Zoom = $source.PageSetup.Zoom
$newSheet.PageSetup.Orientation = $source.PageSetup.Orientation

# Leave the cursor where the author left it when the file was last saved.
[void]$newSheet.Range("E10").Select()
